$d = $word.ActiveDocument

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# 1) Insert the new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaXml = "<w:p $wNs><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Discover the simple gameplay and high volatility of 40 Joker Staxx: 40 Lines. Play for free and find out if this classic slot game is right for you.</w:t></w:r></w:p>"
$metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2) Near the end of the document, drop the duplicated bold title paragraph
#    and rewrite the italic paragraph's text with the new image-prompt copy.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$boldPara = $d.Paragraphs($count - 1)
$italicPara = $d.Paragraphs($count)

$tailRange = $d.Range($boldPara.Range.Start, $italicPara.Range.End)
$featureText = "Create a feature image for `"40 Joker Staxx: 40 Lines`" that highlights the game's modern twist on retro themes. The image should be in a cartoon style, featuring a happy Maya warrior wearing glasses to represent the game's simple yet fun gameplay. The warrior could be holding a classic fruit or a gold ingot, two symbols that represent the game's payout potential. The background could be a mix of retro and modern elements, such as neon lights and classic arcade machines. The overall tone of the image should be vibrant and lively, reflecting the excitement of playing the game."
$tailXml = "<w:p $wNs><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>$featureText</w:t></w:r></w:p>"
$tailRange.InsertXML($tailXml)
